$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 1454.3334
$ws.Range("J2").Value = 1349.25
$ws.Range("L2").Value = 1349.25
$ws.Range("N2").Value = -1575.25

$ws.Range("H40").Value = 2499.889
$ws.Range("I40").Value = 1785.7142
$ws.Range("J40").Value = 4999.5
$ws.Range("K40").Value = 1785.7142
$ws.Range("L40").Value = 4999.5
$ws.Range("M40").Value = -1610.7142
$ws.Range("N40").Value = -5349.5

$ws.Range("H70").Value = 6130.6924
$ws.Range("I70").Value = 2050
$ws.Range("J70").Value = 7944.3335
$ws.Range("K70").Value = 6150
$ws.Range("L70").Value = 23833.0005
$ws.Range("M70").Value = -5880
$ws.Range("N70").Value = -24373.0005

$ws.Range("H73").Value = 6130.6924
$ws.Range("I73").Value = 2050
$ws.Range("J73").Value = 7944.3335
$ws.Range("K73").Value = 6150
$ws.Range("L73").Value = 23833.0005
$ws.Range("M73").Value = -5214
$ws.Range("N73").Value = -25705.0005

$ws.Range("H80").Value = 839.2857
$ws.Range("I80").Value = 887.7778
$ws.Range("J80").Value = 802.9167
$ws.Range("K80").Value = 2663.3334
$ws.Range("L80").Value = 2408.7501
$ws.Range("M80").Value = -1665.3334
$ws.Range("N80").Value = -4404.7501

$ws.Range("H83").Value = 839.2857
$ws.Range("I83").Value = 887.7778
$ws.Range("J83").Value = 802.9167
$ws.Range("K83").Value = 7990.000199999999
$ws.Range("L83").Value = 7226.2503
$ws.Range("M83").Value = -2998.000199999999
$ws.Range("N83").Value = -17210.2503

$ws.Range("H100").Value = 2112.8823
$ws.Range("I100").Value = 2563.3845
$ws.Range("J100").Value = 648.75
$ws.Range("K100").Value = 2563.3845
$ws.Range("L100").Value = 648.75
$ws.Range("M100").Value = -2022.3845
$ws.Range("N100").Value = -1730.75

$ws.Range("H118").Value = 1670.6666
$ws.Range("I118").Value = 576.7143
$ws.Range("J118").Value = 5499.5
$ws.Range("K118").Value = 1730.1429
$ws.Range("L118").Value = 16498.5
$ws.Range("M118").Value = -73.14289999999983
$ws.Range("N118").Value = -19812.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 2147.4583
$ws.Range("I2").Value = 2196.7827
$ws.Range("J2").Value = 1013
$ws.Range("K2").Value = 2196.7827
$ws.Range("L2").Value = 1013
$ws.Range("M2").Value = -2083.7827
$ws.Range("N2").Value = -1239

$ws.Range("H43").Value = 28855.5
$ws.Range("I43").Value = 19617.5
$ws.Range("K43").Value = 19617.5
$ws.Range("M43").Value = -19304.5

$ws.Range("H45").Value = 2270.92
$ws.Range("I45").Value = 1153.1666
$ws.Range("K45").Value = 1153.1666
$ws.Range("M45").Value = -776.1666

$ws.Range("H116").Value = 2147.4583
$ws.Range("I116").Value = 2196.7827
$ws.Range("J116").Value = 1013
$ws.Range("K116").Value = 2196.7827
$ws.Range("L116").Value = 1013
$ws.Range("M116").Value = 97.2172999999998
$ws.Range("N116").Value = -5601

$ws.Range("H122").Value = 1819.45
$ws.Range("I122").Value = 1586.25
$ws.Range("J122").Value = 2752.25
$ws.Range("K122").Value = 4758.75
$ws.Range("L122").Value = 8256.75
$ws.Range("M122").Value = -2308.75
$ws.Range("N122").Value = -13156.75

$ws.Range("H132").Value = 2250.353
$ws.Range("I132").Value = 1219.7778
$ws.Range("K132").Value = 3659.3334
$ws.Range("M132").Value = -1129.3334

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 2147.4583
$ws.Range("I3").Value = 2196.7827
$ws.Range("J3").Value = 1013
$ws.Range("K3").Value = 2196.7827
$ws.Range("L3").Value = 1013
$ws.Range("M3").Value = -2082.7827
$ws.Range("N3").Value = -1241

$ws.Range("H105").Value = 2815.9443
$ws.Range("I105").Value = 2561
$ws.Range("K105").Value = 2561
$ws.Range("M105").Value = -814

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 2011.9412
$ws.Range("I16").Value = 1880.2
$ws.Range("J16").Value = 3000
$ws.Range("K16").Value = 1880.2
$ws.Range("L16").Value = 3000
$ws.Range("M16").Value = -1593.2
$ws.Range("N16").Value = -3574

$ws.Range("H31").Value = 3334573.2
$ws.Range("J31").Value = 1700.6666
$ws.Range("L31").Value = 1700.6666
$ws.Range("N31").Value = -2290.6666

$ws.Range("H34").Value = 3334573.2
$ws.Range("J34").Value = 1700.6666
$ws.Range("L34").Value = 1700.6666
$ws.Range("N34").Value = -2104.6666

$ws.Range("H58").Value = 1828.04
$ws.Range("I58").Value = 1058.4706
$ws.Range("J58").Value = 3463.375
$ws.Range("K58").Value = 1058.4706
$ws.Range("L58").Value = 3463.375
$ws.Range("M58").Value = -855.4706000000001
$ws.Range("N58").Value = -3869.375

$ws.Range("H94").Value = 1896.4736
$ws.Range("J94").Value = 2515.7778
$ws.Range("L94").Value = 2515.7778
$ws.Range("N94").Value = -3417.7778

$ws.Range("H105").Value = 895.1053000000001
$ws.Range("I105").Value = 828.625
$ws.Range("K105").Value = 828.625
$ws.Range("M105").Value = 918.375

$ws.Range("H107").Value = 775.9524
$ws.Range("I107").Value = 516.44446
$ws.Range("K107").Value = 516.44446
$ws.Range("M107").Value = 1403.55554

$ws.Range("H113").Value = 2011.9412
$ws.Range("I113").Value = 1880.2
$ws.Range("J113").Value = 3000
$ws.Range("K113").Value = 1880.2
$ws.Range("L113").Value = 3000
$ws.Range("M113").Value = 289.8
$ws.Range("N113").Value = -7340

$ws.Range("H136").Value = 1828.04
$ws.Range("I136").Value = 1058.4706
$ws.Range("J136").Value = 3463.375
$ws.Range("K136").Value = 3175.4118
$ws.Range("L136").Value = 10390.125
$ws.Range("M136").Value = -625.4118000000003
$ws.Range("N136").Value = -15490.125

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H96").Value = 35255.5
$ws.Range("J96").Value = 35255.5
$ws.Range("L96").Value = 35255.5
$ws.Range("N96").Value = -40747.5

$ws.Range("H113").Value = 2137.6086
$ws.Range("I113").Value = 2480.182
$ws.Range("J113").Value = 1823.5834
$ws.Range("K113").Value = 2480.182
$ws.Range("L113").Value = 1823.5834
$ws.Range("M113").Value = -310.1819999999998
$ws.Range("N113").Value = -6163.5834

$ws.Range("H133").Value = 106000
$ws.Range("J133").Value = 106000
$ws.Range("L133").Value = 106000
$ws.Range("N133").Value = -116120

$ws.Range("H139").Value = 107865.2
$ws.Range("J139").Value = 107865.2
$ws.Range("L139").Value = 107865.2
$ws.Range("N139").Value = -118145.2

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H35").Value = 6547.727
$ws.Range("I35").Value = 4669.4443
$ws.Range("J35").Value = 15000
$ws.Range("K35").Value = 4669.4443
$ws.Range("L35").Value = 15000
$ws.Range("M35").Value = -4333.4443
$ws.Range("N35").Value = -15672

$ws.Range("H43").Value = 6000
$ws.Range("J43").Value = 13333.333
$ws.Range("L43").Value = 13333.333
$ws.Range("N43").Value = -13719.333

$ws.Range("H61").Value = 2678.9285
$ws.Range("I61").Value = 2790.6
$ws.Range("K61").Value = 2790.6
$ws.Range("M61").Value = -2588.6

$ws.Range("H113").Value = 2678.9285
$ws.Range("I113").Value = 2790.6
$ws.Range("K113").Value = 2790.6
$ws.Range("M113").Value = -620.5999999999999

$ws.Range("H136").Value = 3206.5217
$ws.Range("I136").Value = 3252.9285
$ws.Range("J136").Value = 3134.3333
$ws.Range("K136").Value = 9758.7855
$ws.Range("L136").Value = 9402.999899999999
$ws.Range("M136").Value = -7208.7855
$ws.Range("N136").Value = -14502.9999

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H41").Value = 11488.333
$ws.Range("J41").Value = 7757.5
$ws.Range("L41").Value = 7757.5
$ws.Range("N41").Value = -8537.5
